$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.010.24"
$ws.Range("E2").Value = "  +5.58%  "

$ws.Range("D3").Value = "2.459.73"
$ws.Range("E3").Value = "  +3.90%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "573.54"
$ws.Range("E5").Value = "  +2.59%  "

$ws.Range("D6").Value = "145.74"
$ws.Range("E6").Value = "  +6.01%  "

$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("E8").Value = "  +2.80%  "

$ws.Range("D9").Value = "2.458.99"
$ws.Range("E9").Value = "  +4.07%  "

$ws.Range("E10").Value = "  +6.14%  "

$ws.Range("D11").Value = "0.161"
$ws.Range("E11").Value = "  +1.10%  "

$ws.Range("E12").Value = "  +3.08%  "

$ws.Range("E13").Value = "  +5.18%  "

$ws.Range("D14").Value = "27.38"
$ws.Range("E14").Value = "  +7.51%  "

$ws.Range("E15").Value = "  +8.25%  "

$ws.Range("D16").Value = "2.905.29"
$ws.Range("E16").Value = "  +3.97%  "

$ws.Range("D17").Value = "62.893.77"
$ws.Range("E17").Value = "  +5.31%  "

$ws.Range("D18").Value = "2.458.43"
$ws.Range("E18").Value = "  +3.61%  "

$ws.Range("D19").Value = "7.88"
$ws.Range("E19").Value = "  -1.90%  "

$ws.Range("E20").Value = "  +5.24%  "

$ws.Range("D21").Value = "328.60"
$ws.Range("E21").Value = "  +2.22%  "

$ws.Range("E22").Value = "  +2.39%  "

$ws.Range("E23").Value = "  +12.86%  "

$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.19%  "

$ws.Range("E25").Value = "  +2.56%  "

$ws.Range("D26").Value = "634.46"
$ws.Range("E26").Value = "  +13.62%  "

$ws.Range("E27").Value = "  +10.36%  "

$ws.Range("D28").Value = "8.49"
$ws.Range("E28").Value = "  +4.49%  "

$ws.Range("D29").Value = "0.0₃0984"
$ws.Range("E29").Value = "  +7.09%  "

$ws.Range("D30").Value = "2.533.01"

$ws.Range("E31").Value = "  +2.80%  "

$ws.Range("E32").Value = "  +9.08%  "

$ws.Range("E33").Value = "  +5.92%  "

$ws.Range("D34").Value = "1.85"
$ws.Range("E34").Value = "  +4.62%  "

$ws.Range("E35").Value = "  +5.13%  "

$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("E37").Value = "  +5.03%  "

$ws.Range("D38").Value = "0.374"
$ws.Range("E38").Value = "  +2.30%  "

$ws.Range("D39").Value = "152.07"
$ws.Range("E39").Value = "  -0.58%  "

$ws.Range("D40").Value = "5.40"
$ws.Range("E40").Value = "  +8.90%  "

$ws.Range("D41").Value = "18.70"
$ws.Range("E41").Value = "  +3.25%  "

$ws.Range("E42").Value = "  +14.29%  "

$ws.Range("E43").Value = "  +7.80%  "

$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("D45").Value = "0.0₆0296"
$ws.Range("E45").Value = "  -0.50%  "

$ws.Range("D46").Value = "144.96"
$ws.Range("E46").Value = "  +4.86%  "

$ws.Range("E47").Value = "  +2.41%  "

$ws.Range("D48").Value = "20.35"
$ws.Range("E48").Value = "  +6.91%  "

$ws.Range("D49").Value = "0.603"
$ws.Range("E49").Value = "  +3.22%  "

$ws.Range("D50").Value = "0.0517"
$ws.Range("E50").Value = "  +3.59%  "

$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "12.45"
$ws.Range("E51").Value = "  +6.69%  "
